# Appends the 13 match rows that were re-scraped/re-appended below the
# existing data (rows 2-14 duplicated, in shuffled order, as rows 15-27),
# matching the author's "json file code updated" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

# Every value in this sheet is stored as text (see the existing
# numberStoredAsText ignoredError on A1:K14), including the numeric-
# looking columns (runs, balls, 4s, 6s, strike rate) - e.g. "140.00"
# must stay "140.00", not become the number 140. Setting the cell
# NumberFormat to Text ("@") before assigning the value forces that.
$newRange = $ws.Range("A15:K27")
$newRange.NumberFormat = "@"

$rows = @(
    @(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Mumbai Indians", "Rajasthan Royals", "Hardik Pandya ", "60", "21", "2", "7", "285.71"),
    @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Mumbai Indians", "Chennai Super Kings", "Hardik Pandya ", "14", "10", "0", "2", "140.00"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Mumbai Indians", "Royal Challengers Bangalore", "Hardik Pandya ", "17", "15", "0", "2", "113.33"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Mumbai Indians", "Kolkata Knight Riders", "Hardik Pandya ", "18", "13", "2", "1", "138.46"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Hardik Pandya ", "0", "2", "0", "0", "0.00"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Delhi Capitals", "Hardik Pandya ", "37", "14", "0", "5", "264.28"),
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Mumbai Indians", "Kings XI Punjab", "Hardik Pandya ", "30", "11", "3", "2", "272.72"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Mumbai Indians", "Sunrisers Hyderabad", "Hardik Pandya ", "28", "19", "2", "2", "147.36"),
    @(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Mumbai Indians", "Kolkata Knight Riders", "Hardik Pandya ", "21", "11", "3", "1", "190.90"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Hardik Pandya ", "3", "5", "0", "0", "60.00"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Mumbai Indians", "Royal Challengers Bangalore", "Hardik Pandya ", "15", "13", "0", "1", "115.38"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Rajasthan Royals", "Hardik Pandya ", "30", "19", "2", "1", "157.89"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Mumbai Indians", "Kings XI Punjab", "Hardik Pandya ", "8", "4", "0", "1", "200.00")
)

$startRow = 15
for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $rows[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $ws.Range($columns[$j] + $rowNum).Value = $rowValues[$j]
    }
}

